# The fixture ("Proximos" / upcoming-matches) was re-scraped: the whole
# table shifted up by one row (an obsolete blank spacer row was removed)
# and the now-stale last match row (week 38, 2026-05-24 Man City vs
# Aston Villa) fell off the bottom of the sheet.
#
# Deleting any one of the always-blank rows above the first fixture row
# reproduces that shift exactly (every row from the deleted one down
# moves up by one, the sheet's dimension shrinks from L159 to L158, and
# the very last row disappears) without touching any cell's stored
# value/type - so dates stored as plain text stay plain text instead of
# being re-interpreted as Excel date serials.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Delete()
